# Add a new "interesting paper" link (arxiv) right after the existing
# ResearchGate hyperlink paragraph, separated by one blank paragraph,
# matching the author's commit "Add files via upload / Interesting papers".

$d = $word.ActiveDocument

# Locate the paragraph that holds the ResearchGate hyperlink so the new
# content is anchored robustly (not by a hard-coded paragraph index).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*researchgate.net/publication/324957692_Financial_Sentiment_Lexicon_Analysis*") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not find the ResearchGate hyperlink paragraph to anchor the new content."
}

$anchorIndex = $anchor.Index

# Insert a blank paragraph right after the ResearchGate paragraph.
$anchor.Range.InsertParagraphAfter()

# Insert another new paragraph after that blank one - this will hold the
# new arxiv hyperlink.
$blankPara = $d.Paragraphs.Item($anchorIndex + 1)
$blankPara.Range.InsertParagraphAfter()

$linkPara = $d.Paragraphs.Item($anchorIndex + 2)

$url = "https://arxiv.org/pdf/1307.5336.pdf"

# Put the URL text into the new paragraph, then convert just that text
# into a hyperlink run (mirrors how the other hyperlinks in the doc are
# built: a single run styled with the built-in "Hyperlink" character style).
$linkPara.Range.InsertAfter($url)
$linkRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start + $url.Length)
$d.Hyperlinks.Add($linkRange, $url) | Out-Null

Write-Output "Inserted arxiv hyperlink paragraph after paragraph $anchorIndex"
